$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in column H, matching the style of the other headers (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Populate the "Save" indicator column: 1 when the "sum" (column G) is >= 8, else 0
for ($r = 2; $r -le 75; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -ge 8) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
